# US 3.3 commit files
# Applies the changes described by the diff to the workbook:
#  - "About" sheet: add Notes section (rows 9-11)
#  - "PPEIdtICEaT" sheet: update header cells A1/B1 (text + wrap formatting)
#  - "Data" sheet: no content changes required (only incidental style/string
#    re-indexing happens automatically as a result of the edits elsewhere)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PPEIdtICEaT sheet: update the table header
# ---------------------------------------------------------------------
$wsPPE = $wb.Worksheets.Item("PPEIdtICEaT")

# A1 header becomes the fuller "Building Component Efficiency Improvement
# (dimensionless)" label, and both header cells wrap their text so the
# header row grows taller.
$wsPPE.Range("A1").Value = "Building Component Efficiency Improvement (dimensionless)"
$wsPPE.Range("B1").Value = "Efficiency Improvement"

$wsPPE.Range("A1").WrapText = $true
$wsPPE.Range("B1").WrapText = $true

$wsPPE.Rows.Item(1).RowHeight = 45

# ---------------------------------------------------------------------
# About sheet: append a Notes section explaining which building
# component the policy maps to
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A1").Value = "PPEIdtICEaT Potential Percentage Eff Improvement due to Improved Contractor Edu and Training"

$wsAbout.Range("A9").Value = "Notes:"
$wsAbout.Range("A10").Value = "This policy covers improvements in air sealing, framing, and insulation, so it applies to the "
$wsAbout.Range("A11").Value = """envelope"" component."

$wsAbout.Activate()
$wsAbout.Range("A11").Select()
